# Add the 250-packet run data (new "230 Packets" labeled block in columns M:O / Q:S)
# on the "MW CNN 1D" worksheet (sheet3), per commit "added 250 packet run MW 1D".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the data block (row 1) and column headers (row 2)
$ws.Range("R1").Value = "230 Packets"
$ws.Range("Q2").Value = "Acc"
$ws.Range("R2").Value = "Loss"
$ws.Range("S2").Value = "Time"

# Per-trial Acc / Loss / Time values for the new run, columns M:O, rows 3-51
$newRunData = @(
    @(3, 88.483887910842896, 0.281325041714515, 551.36055135726895),
    @(4, 88.996136188506995, 0.252420932368148, 544.04951500892605),
    @(5, 88.417482376098604, 0.306426730454052, 542.28713583946205),
    @(6, 89.041191339492798, 0.26398389835938002, 546.41168141364994),
    @(7, 88.889414072036701, 0.28471131582523101, 541.61022663116398),
    @(8, 88.851469755172701, 0.27909826622030498, 542.92116403579701),
    @(9, 88.986647129058795, 0.26166537757056002, 545.85422968864395),
    @(10, 88.355821371078406, 0.284870725911266, 547.05730724334705),
    @(11, 88.597720861434894, 0.28846056012302701, 544.99331116676296),
    @(12, 88.6830925941467, 0.26560601756776298, 546.16639089584305),
    @(13, 87.36452460289, 0.359610378499442, 546.45408487319901),
    @(14, 88.901275396347003, 0.26030729608185499, 546.16273212432804),
    @(15, 88.614320755004798, 0.28436701727602598, 546.87314581870999),
    @(16, 89.036452770233097, 0.27837585723352398, 547.32949757575898),
    @(17, 88.939219713211003, 0.26104552296434202, 545.30119872093201),
    @(18, 88.711553812026906, 0.26672030669906699, 545.76163697242703),
    @(19, 88.640403747558594, 0.25990669919340398, 545.28115653991699),
    @(20, 89.010363817214895, 0.27289671628287598, 545.80812358856201),
    @(21, 88.585859537124605, 0.29010992726569801, 547.65209531784001),
    @(22, 88.204044103622394, 0.30982169429669798, 545.26629543304398),
    @(23, 87.639623880386296, 0.264758731574678, 544.52358675002995),
    @(24, 90.044349431991506, 0.24955440453571301, 546.63705778121903),
    @(25, 89.543956518173204, 0.25332823623609002, 546.20874142646699),
    @(26, 88.834869861602698, 0.26073585859903697, 547.71153569221497),
    @(27, 82.941639423370304, 2.6234215325207999, 545.36057758331299),
    @(28, 82.941639423370304, 2.62342143300845, 543.04964137077297),
    @(29, 88.777953386306706, 0.27444691853879799, 544.50223922729401),
    @(30, 88.808786869049001, 0.28404082121221902, 547.74136209487904),
    @(31, 88.626176118850694, 0.26247188221237999, 543.99815678596497),
    @(32, 89.100480079650794, 0.259553523475755, 547.98098540306),
    @(33, 88.846731185913001, 0.26531235577315598, 544.29760575294495),
    @(34, 88.633292913436804, 0.27672774358729701, 544.23738980293194),
    @(35, 89.067280292510901, 0.248404741251362, 543.86246681213299),
    @(36, 89.031708240509005, 0.26492855713896202, 544.77868914604096),
    @(37, 88.830125331878605, 0.25159506870326198, 547.13121771812405),
    @(38, 87.608796358108506, 0.29609051382293999, 547.16282296180702),
    @(39, 87.947922945022498, 0.28873343550453801, 549.94790124893098),
    @(40, 88.426971435546804, 0.27449781699912401, 545.10195541381802),
    @(41, 89.043563604354802, 0.268049381352639, 547.15556144714299),
    @(42, 88.846731185913001, 0.25549514489000302, 552.29076719284001),
    @(43, 88.962930440902696, 0.26159132634100402, 547.493482589721),
    @(44, 88.661748170852604, 0.28134602580025903, 546.67198348045304),
    @(45, 88.936847448348999, 0.24437450679058301, 542.47305011749199),
    @(46, 89.088624715804997, 0.25561373741417698, 547.44016408920197),
    @(47, 89.000874757766695, 0.26804574518342, 544.81893038749695),
    @(48, 88.901275396347003, 0.25049732122971802, 546.60061383247296),
    @(49, 88.974791765212998, 0.26366006618967502, 546.98892784118596),
    @(50, 88.657009601593003, 0.27222664157325699, 543.56502509117104),
    @(51, 89.638817310333195, 0.25382578302136999, 546.291284561157)
)

foreach ($entry in $newRunData) {
    $row = $entry[0]
    $ws.Cells.Item($row, 13).Value = $entry[1]
    $ws.Cells.Item($row, 14).Value = $entry[2]
    $ws.Cells.Item($row, 15).Value = $entry[3]
}

# Match the saved selection state recorded in the workbook
$ws.Range("T12").Select()
